$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 134
$ws.Range("H134").Value = 45973.777
$ws.Range("J134").Value = 45973.777
$ws.Range("L134").Value = 45973.777
$ws.Range("N134").Value = -56113.777

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3236.4167
$ws.Range("I2").Value = 1716.7222
$ws.Range("J2").Value = 7795.5
$ws.Range("K2").Value = 1716.7222
$ws.Range("L2").Value = 7795.5
$ws.Range("M2").Value = -1603.7222
$ws.Range("N2").Value = -8021.5

# Row 3
$ws.Range("H3").Value = 4500
$ws.Range("I3").Value = 3666.6667
$ws.Range("K3").Value = 3666.6667
$ws.Range("M3").Value = -3551.6667

# Row 7
$ws.Range("H7").Value = 29282.908
$ws.Range("J7").Value = 29282.908
$ws.Range("L7").Value = 29282.908
$ws.Range("N7").Value = -29510.908

# Row 45
$ws.Range("H45").Value = 3342.8845
$ws.Range("I45").Value = 2169.8823
$ws.Range("J45").Value = 5558.5557
$ws.Range("K45").Value = 2169.8823
$ws.Range("L45").Value = 5558.5557
$ws.Range("M45").Value = -1792.8823
$ws.Range("N45").Value = -6312.5557

# Row 52
$ws.Range("H52").Value = 14733.333
$ws.Range("J52").Value = 14733.333
$ws.Range("L52").Value = 14733.333
$ws.Range("N52").Value = -15369.333

# Row 74
$ws.Range("H74").Value = 1059.6
$ws.Range("I74").Value = 1068.2307
$ws.Range("J74").Value = 1003.5
$ws.Range("K74").Value = 1068.2307
$ws.Range("L74").Value = 1003.5
$ws.Range("M74").Value = -194.2307000000001
$ws.Range("N74").Value = -2751.5

# Row 77
$ws.Range("H77").Value = 1059.6
$ws.Range("I77").Value = 1068.2307
$ws.Range("J77").Value = 1003.5
$ws.Range("K77").Value = 5341.1535
$ws.Range("L77").Value = 5017.5
$ws.Range("M77").Value = -973.1535000000003
$ws.Range("N77").Value = -13753.5

# Row 116
$ws.Range("H116").Value = 3236.4167
$ws.Range("I116").Value = 1716.7222
$ws.Range("J116").Value = 7795.5
$ws.Range("K116").Value = 1716.7222
$ws.Range("L116").Value = 7795.5
$ws.Range("M116").Value = 577.2778000000001
$ws.Range("N116").Value = -12383.5

# Row 132
$ws.Range("H132").Value = 84578.38
$ws.Range("I132").Value = 101098.98
$ws.Range("K132").Value = 303296.94
$ws.Range("M132").Value = -300766.94

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3236.4167
$ws.Range("I3").Value = 1716.7222
$ws.Range("J3").Value = 7795.5
$ws.Range("K3").Value = 1716.7222
$ws.Range("L3").Value = 7795.5
$ws.Range("M3").Value = -1602.7222
$ws.Range("N3").Value = -8023.5

# Row 8
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = ""

# Row 56
$ws.Range("H56").Value = 30110
$ws.Range("J56").Value = 30110
$ws.Range("L56").Value = 30110
$ws.Range("N56").Value = -31588

# Row 102
$ws.Range("H102").Value = 25000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 25000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 25000
$ws.Range("M102").Value = ""
$ws.Range("N102").Value = -31490

# Row 128
$ws.Range("H128").Value = 1613.3334
$ws.Range("I128").Value = 1613.3334
$ws.Range("K128").Value = 4840.0002
$ws.Range("M128").Value = -2350.0002

# Row 134
$ws.Range("H134").Value = 112669.7
$ws.Range("I134").Value = 137608.27
$ws.Range("J134").Value = 2940
$ws.Range("K134").Value = 412824.8099999999
$ws.Range("L134").Value = 8820
$ws.Range("M134").Value = -410289.8099999999
$ws.Range("N134").Value = -13890

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 2333.3333
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -2887
$ws.Range("N3").Value = -2226

# Row 31
$ws.Range("H31").Value = 8334836.5
$ws.Range("I31").Value = 1275.6923
$ws.Range("J31").Value = 23811450
$ws.Range("K31").Value = 1275.6923
$ws.Range("L31").Value = 23811450
$ws.Range("M31").Value = -980.6922999999999
$ws.Range("N31").Value = -23812040

# Row 34
$ws.Range("H34").Value = 8334836.5
$ws.Range("I34").Value = 1275.6923
$ws.Range("J34").Value = 23811450
$ws.Range("K34").Value = 1275.6923
$ws.Range("L34").Value = 23811450
$ws.Range("M34").Value = -1073.6923
$ws.Range("N34").Value = -23811854

$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Range("H44").Value = 1739.8572
$ws.Range("J44").Value = 1739.8572
$ws.Range("L44").Value = 5219.571599999999
$ws.Range("N44").Value = -6015.571599999999

# Row 123
$ws.Range("H123").Value = 4629.8335
$ws.Range("I123").Value = 1243.3334
$ws.Range("J123").Value = 5113.619
$ws.Range("K123").Value = 3730.0002
$ws.Range("L123").Value = 15340.857
$ws.Range("M123").Value = -1280.0002
$ws.Range("N123").Value = -20240.857

$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Range("H32").Value = 32270
$ws.Range("J32").Value = 32270
$ws.Range("L32").Value = 32270
$ws.Range("N32").Value = -32862

# Row 59
$ws.Range("H59").Value = 1687.5
$ws.Range("I59").Value = 1687.5
$ws.Range("K59").Value = 1687.5
$ws.Range("M59").Value = -1104.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5002404
$ws.Range("I7").Value = 6668645.5
$ws.Range("J7").Value = 3679.8
$ws.Range("K7").Value = 6668645.5
$ws.Range("L7").Value = 3679.8
$ws.Range("M7").Value = -6668533.5
$ws.Range("N7").Value = -3903.8

# Row 126
$ws.Range("H126").Value = 5002404
$ws.Range("I126").Value = 6668645.5
$ws.Range("J126").Value = 3679.8
$ws.Range("K126").Value = 20005936.5
$ws.Range("L126").Value = 11039.4
$ws.Range("M126").Value = -20003466.5
$ws.Range("N126").Value = -15979.4

# Row 136
$ws.Range("H136").Value = 1836.2593
$ws.Range("I136").Value = 1617.4546
$ws.Range("J136").Value = 2799
$ws.Range("K136").Value = 4852.3638
$ws.Range("L136").Value = 8397
$ws.Range("M136").Value = -2302.3638
$ws.Range("N136").Value = -13497

$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 3769476
$ws.Range("I3").Value = 15000000
$ws.Range("J3").Value = 25968
$ws.Range("K3").Value = 15000000
$ws.Range("L3").Value = 25968
$ws.Range("M3").Value = -14999886
$ws.Range("N3").Value = -26196

# Row 11
$ws.Range("H11").Value = 49336.668
$ws.Range("I11").Value = 8000
$ws.Range("J11").Value = 70005
$ws.Range("K11").Value = 8000
$ws.Range("L11").Value = 70005
$ws.Range("M11").Value = -7858
$ws.Range("N11").Value = -70289

# Row 119
$ws.Range("H119").Value = 23499
$ws.Range("J119").Value = 23499
$ws.Range("L119").Value = 23499
$ws.Range("N119").Value = -33175

# Row 126
$ws.Range("H126").Value = 1700.8
$ws.Range("I126").Value = 1523.409
$ws.Range("J126").Value = 3001.6667
$ws.Range("K126").Value = 4570.227000000001
$ws.Range("L126").Value = 9005.000100000001
$ws.Range("M126").Value = -2100.227000000001
$ws.Range("N126").Value = -13945.0001
